$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (this also updates the _FilterDatabase defined name reference)
$ws.Name = "2019_As"

# Update the pollutant value (C2) to the new pollutant "As"
$ws.Range("C2").Value = "As"

# Update the selection to match the authored state
$ws.Range("C3").Select()
